# Single-threaded Excel update: slide the rolling data window up by one
# row (dropping the oldest row, row 5) and append the newest day's data
# into row 9. Also restores the active selection that Excel records when
# a user highlights the moved block (A5:N8) before saving.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the current values of rows 6-9 (columns A:M) before overwriting
# rows 5-8, since row 5 will be overwritten first. Value2 is used because
# it round-trips plain numbers reliably (Value can box them oddly).
$numCols = 13  # columns A.. M
$snapshot = @{}
for ($r = 6; $r -le 9; $r++) {
    $rowVals = @()
    for ($c = 1; $c -le $numCols; $c++) {
        $rowVals += , $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

# Shift rows 6-9 up into rows 5-8.
for ($r = 6; $r -le 9; $r++) {
    $destRow = $r - 1
    $rowVals = $snapshot[$r]
    for ($c = 1; $c -le $numCols; $c++) {
        $ws.Cells.Item($destRow, $c).Value2 = $rowVals[$c - 1]
    }
}

# Write the newest day's data into row 9 (column N keeps its existing
# "Bag" label so it doesn't need to be touched).
$newRow = @(42612.889490740738, 10, 50, 45, 74, 25, 11470, 21795, 2634, 297, 266, 35, 12)
for ($c = 1; $c -le $numCols; $c++) {
    $ws.Cells.Item(9, $c).Value2 = $newRow[$c - 1]
}

# Record the selection left behind from highlighting the shifted block.
$ws.Range("A5:N8").Select()
